$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing content (keep formatting) then drop the extra row/column
# that the new data no longer needs.
$ws.UsedRange.ClearContents()
$ws.Rows.Item(3).Delete()
$ws.Columns.Item(25).Delete()

# New header row (row 1) — metrics reordered (time, rew, waitingTime) and a
# rule-based-agent testing column ("O_rule_set") now populated.
$headers = @(
    "time","rew","waitingTime",
    "E_num_train_rollouts","E_rollout_length","E_eval_freq","E_eval_num_eps","E_max_ep_steps","E_test_num_eps",
    "M_state_size","M_action_size",
    "A_gae_tau","A_entropy_weight","A_minibatch_size","A_optimization_epochs","A_ppo_ratio_clip","A_discount","A_learning_rate","A_clip_grads","A_gradient_clip","A_value_loss_coef",
    "O_num_agents","O_rule_set","O_rule_set_params"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# New data row (row 2) for the rule-based agent test run.
$values = @(
    0.1966046134630839, 870.8666666666667, 12.69385964912281,
    20, 64, 5, 5, 250, 20,
    6, 2,
    0.9, 0.01, 32, 10, 0.3, 0.99, 0.001, $true, 1, 1,
    3, "timer", "{'length': 10}"
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}
